$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 51.714287
$ws.Range("I5").Value = 64
$ws.Range("K5").Value = 64
$ws.Range("M5").Value = 51

$ws.Range("H43").Value = 8500
$ws.Range("I43").Value = 7750
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 7750
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -7681
$ws.Range("N43").Value = -10138

$ws.Range("H48").Value = 10000
$ws.Range("J48").Value = 10000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30584

$ws.Range("H56").Value = 10000
$ws.Range("J56").Value = 10000
$ws.Range("L56").Value = 30000
$ws.Range("N56").Value = -31068

$ws.Range("H86").Value = 4862.25
$ws.Range("J86").Value = 4899.6665
$ws.Range("L86").Value = 4899.6665
$ws.Range("N86").Value = -7145.6665

$ws.Range("H89").Value = 4862.25
$ws.Range("J89").Value = 4899.6665
$ws.Range("L89").Value = 24498.3325
$ws.Range("N89").Value = -35730.3325

$ws.Range("H98").Value = 2325.375
$ws.Range("I98").Value = 2325.375
$ws.Range("K98").Value = 2325.375
$ws.Range("M98").Value = -827.375

$ws.Range("H116").Value = 6610.75
$ws.Range("I116").Value = 6648
$ws.Range("J116").Value = 6499
$ws.Range("K116").Value = 6648
$ws.Range("L116").Value = 6499
$ws.Range("M116").Value = -3206
$ws.Range("N116").Value = -13383

$ws.Range("H122").Value = 2325.375
$ws.Range("I122").Value = 2325.375
$ws.Range("K122").Value = 6976.125
$ws.Range("M122").Value = -4526.125

$ws.Range("H132").Value = 2787.0667
$ws.Range("I132").Value = 2446.6155
$ws.Range("K132").Value = 7339.8465
$ws.Range("M132").Value = -4809.8465

$ws.Range("H138").Value = 4546.857
$ws.Range("I138").Value = 1221.875
$ws.Range("J138").Value = 5876.85
$ws.Range("K138").Value = 3665.625
$ws.Range("L138").Value = 17630.55
$ws.Range("M138").Value = 1474.375
$ws.Range("N138").Value = -27910.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5256.1704
$ws.Range("I32").Value = 3654.1765
$ws.Range("J32").Value = 9446
$ws.Range("K32").Value = 3654.1765
$ws.Range("L32").Value = 9446
$ws.Range("M32").Value = -3367.1765
$ws.Range("N32").Value = -10020

$ws.Range("H122").Value = 359929.7
$ws.Range("I122").Value = 529062.5
$ws.Range("K122").Value = 1587187.5
$ws.Range("M122").Value = -1584737.5

$ws.Range("H132").Value = 2039.5
$ws.Range("I132").Value = 2028.4286
$ws.Range("J132").Value = 2065.3333
$ws.Range("K132").Value = 6085.2858
$ws.Range("L132").Value = 6195.999899999999
$ws.Range("M132").Value = -3555.2858
$ws.Range("N132").Value = -11255.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1434.4
$ws.Range("I99").Value = 1296.037
$ws.Range("K99").Value = 1296.037
$ws.Range("M99").Value = 201.963

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 883
$ws.Range("I70").Value = 883
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2649
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -2334

$ws.Range("H73").Value = 883
$ws.Range("I73").Value = 883
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2649
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -1557

$ws.Range("H113").Value = 3199.75
$ws.Range("I113").Value = 4899
$ws.Range("J113").Value = 2633.3333
$ws.Range("K113").Value = 14697
$ws.Range("L113").Value = 7899.999899999999
$ws.Range("M113").Value = -12527
$ws.Range("N113").Value = -12239.9999

$ws.Range("H120").Value = 12696.286
$ws.Range("I120").Value = 8291.333000000001
$ws.Range("K120").Value = 24873.999
$ws.Range("M120").Value = -20035.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 465.5
$ws.Range("I2").Value = 72.71429000000001
$ws.Range("J2").Value = 771
$ws.Range("K2").Value = 72.71429000000001
$ws.Range("L2").Value = 771
$ws.Range("M2").Value = 40.28570999999999
$ws.Range("N2").Value = -997

$ws.Range("H96").Value = 54951
$ws.Range("J96").Value = 54951
$ws.Range("L96").Value = 54951
$ws.Range("N96").Value = -60443

$ws.Range("H122").Value = 79887.46000000001
$ws.Range("I122").Value = 2873.8
$ws.Range("K122").Value = 8621.400000000001
$ws.Range("M122").Value = -6171.400000000001

$ws.Range("H123").Value = 26850.154
$ws.Range("J123").Value = 26850.154
$ws.Range("L123").Value = 26850.154
$ws.Range("N123").Value = -31750.154

$ws.Range("H132").Value = 2130.162
$ws.Range("I132").Value = 1709.4445
$ws.Range("J132").Value = 2528.7368
$ws.Range("K132").Value = 5128.333500000001
$ws.Range("L132").Value = 7586.2104
$ws.Range("M132").Value = -2598.333500000001
$ws.Range("N132").Value = -12646.2104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4498.75
$ws.Range("I7").Value = 4498.75
$ws.Range("K7").Value = 4498.75
$ws.Range("M7").Value = -4386.75

$ws.Range("H55").Value = 464.66666
$ws.Range("I55").Value = 464.66666
$ws.Range("K55").Value = 464.66666
$ws.Range("M55").Value = -291.66666

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 4544.5835
$ws.Range("I122").Value = 4318.5557
$ws.Range("J122").Value = 5222.6665
$ws.Range("K122").Value = 12955.6671
$ws.Range("L122").Value = 15667.9995
$ws.Range("M122").Value = -10505.6671
$ws.Range("N122").Value = -20567.9995

$ws.Range("H126").Value = 4498.75
$ws.Range("I126").Value = 4498.75
$ws.Range("K126").Value = 13496.25
$ws.Range("M126").Value = -11026.25

$ws.Range("H136").Value = 4386.778
$ws.Range("J136").Value = 4717.8
$ws.Range("L136").Value = 14153.4
$ws.Range("N136").Value = -19253.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2261.6667
$ws.Range("I126").Value = 1708.9286
$ws.Range("K126").Value = 5126.7858
$ws.Range("M126").Value = -2656.7858

$ws.Range("H132").Value = 31473.451
$ws.Range("I132").Value = 38651.72
$ws.Range("K132").Value = 115955.16
$ws.Range("M132").Value = -113425.16
